# The commit swaps the two embedded theme parts (ppt/theme/theme1.xml,
# which held the "Office Theme" color scheme, and ppt/theme/theme2.xml,
# which held the "Integral" color scheme) so that the deck's live theme
# (the one bound to the slide master / all slides) now carries the
# "Office Theme" palette instead of "Integral".
#
# The PowerPoint object model doesn't expose a way to rename the raw
# <a:theme>/<a:clrScheme> name attributes, but it does let us rewrite
# every one of the 12 theme colors (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) through Slide.ThemeColorScheme, which is exactly the
# substantive, visible part of the swap for the theme that is actually
# in use across the deck.

$p = $ppt.ActivePresentation

# Target palette = the "Office Theme" colors that used to live in
# ppt/theme/theme1.xml, expressed as RRGGBB hex, in theme-color order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 1; $i -le $tcs.Count; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $oleColor = ($b * 65536) + ($g * 256) + $r

    $tcs.Colors($i).RGB = $oleColor
}
